$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) cells are stored as plain text in the workbook
# (e.g. "1.00", "0.999", "65.032.71"). Excel auto-detects numeric-looking
# strings and would silently convert them to numbers (losing the original
# textual formatting) unless the cell is pre-formatted as Text. Apply the
# Text format only to the specific cells being updated.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.948.90"
$ws.Range("E2").Value = "  -2.67%  "
$ws.Range("D3").Value = "3.184.73"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "603.05"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").Value = "152.90"
$ws.Range("E6").Value = "  -4.33%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "3.183.22"
$ws.Range("E8").Value = "  -1.66%  "
$ws.Range("E9").Value = "  -3.29%  "
$ws.Range("E10").Value = "  -4.71%  "
$ws.Range("D11").Value = "5.61"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("D12").Value = "0.480"
$ws.Range("E12").Value = "  -5.48%  "
$ws.Range("D13").Value = "0.0000263"
$ws.Range("E13").Value = "  -3.94%  "
$ws.Range("D14").Value = "37.29"
$ws.Range("E14").Value = "  -4.80%  "
$ws.Range("D15").Value = "3.689.22"
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("D16").Value = "65.004.07"
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("D17").Value = "3.187.95"
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "7.10"
$ws.Range("E19").Value = "  -4.18%  "
$ws.Range("D20").Value = "485.72"
$ws.Range("E20").Value = "  -5.23%  "
$ws.Range("D21").Value = "15.00"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").Value = "0.722"
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("D23").Value = "7.85"
$ws.Range("E23").Value = "  -2.53%  "
$ws.Range("D24").Value = "14.13"
$ws.Range("E24").Value = "  -3.99%  "
$ws.Range("D25").Value = "85.21"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").Value = "2.95"
$ws.Range("E27").Value = "  -2.12%  "
$ws.Range("D28").Value = "8.84"
$ws.Range("E28").Value = "  -3.37%  "
$ws.Range("D29").Value = "2.27"
$ws.Range("E29").Value = "  -4.45%  "
$ws.Range("D30").Value = "7.31"
$ws.Range("E30").Value = "  +3.83%  "
$ws.Range("E31").Value = "  +1.85%  "
$ws.Range("D32").Value = "2.74"
$ws.Range("E32").Value = "  -7.56%  "
$ws.Range("D33").Value = "27.06"
$ws.Range("E33").Value = "  -4.38%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  -5.42%  "
$ws.Range("D36").Value = "6.18"
$ws.Range("E36").Value = "  -5.48%  "
$ws.Range("D37").Value = "3.31"
$ws.Range("E37").Value = "  +7.34%  "
$ws.Range("D38").Value = "54.64"
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("D39").Value = "0.0₃0751"
$ws.Range("E39").Value = "  -3.13%  "
$ws.Range("D40").Value = "463.14"
$ws.Range("E40").Value = "  -9.21%  "
$ws.Range("D41").Value = "0.128"
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("D42").Value = "0.0408"
$ws.Range("E42").Value = "  -3.60%  "
$ws.Range("D43").Value = "8.59"
$ws.Range("E43").Value = "  -2.24%  "
$ws.Range("D44").Value = "2.47"
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("D45").Value = "2.934.35"
$ws.Range("E45").Value = "  +0.78%  "
$ws.Range("D46").Value = "0.279"
$ws.Range("E46").Value = "  -7.12%  "
$ws.Range("D47").Value = "27.44"
$ws.Range("E47").Value = "  -3.25%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").Value = "2.37"
$ws.Range("E49").Value = "  -3.14%  "
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").Value = "120.90"
$ws.Range("E51").Value = "  -2.03%  "
